$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.238.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.860.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.93%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'242.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.6986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.07830"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.20%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3123"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.67%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'24.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.71%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -4.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.865.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.137"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'91.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.73%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.88%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +3.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008516"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.37%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'29.275.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.43%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'248.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.53%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.115.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -3.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.579"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.95%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.08%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1541"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.68%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'160.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.86%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.911"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.79%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.578"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.66%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.95%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.242"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.25%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.207"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7618"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.50%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.881"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.31%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E37").Value = "'  -0.16%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.98%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.244.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.76%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.9016"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.36%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'110.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.908"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -7.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = "'0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'68.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.012.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -3.85%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.567"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.60%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.5182"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.14%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.93%  "
$ws.Range("E51").Style = "Normal"
